# fix shark double counts
$wb = $excel.ActiveWorkbook

# --- Sheet: "Status by Landings (Area)" ---
$wsArea = $wb.Worksheets.Item("Status by Landings (Area)")
$wsArea.Range("C2").Value = 1.124014212475507
$wsArea.Range("C5").Value = 3.028873523761085
$wsArea.Range("C7").Value = 25.99140387354132
$wsArea.Range("C8").Value = 44.04745699154409
$wsArea.Range("C9").Value = 29.9611391349146
$wsArea.Range("C10").Value = 70.03886086508541
$wsArea.Range("C11").Value = 29.9611391349146

# --- Sheet: "Status by Landings (Tier)" ---
$wsTier = $wb.Worksheets.Item("Status by Landings (Tier)")

$wsTier.Range("B4").Value = 1.124014212475507
$wsTier.Range("E4").Value = 3.028873523761085
$wsTier.Range("G4").Value = 25.99140387354132
$wsTier.Range("H4").Value = 44.04745699154409
$wsTier.Range("I4").Value = 29.9611391349146
$wsTier.Range("J4").Value = 70.03886086508541
$wsTier.Range("K4").Value = 29.9611391349146

$wsTier.Range("B5").Value = 1.124014212475507
$wsTier.Range("E5").Value = 3.028873523761085
$wsTier.Range("G5").Value = 25.99140387354132
$wsTier.Range("H5").Value = 44.04745699154409
$wsTier.Range("I5").Value = 29.9611391349146
$wsTier.Range("J5").Value = 70.03886086508541
$wsTier.Range("K5").Value = 29.9611391349146

# --- Sheet: "Comparison by Landings" ---
$wsComp = $wb.Worksheets.Item("Comparison by Landings")
$wsComp.Range("C2").Value = 95.35980609945506
$wsComp.Range("C3").Value = 25.99140387354132
$wsComp.Range("C4").Value = 44.04745699154409
$wsComp.Range("C5").Value = 29.9611391349146
$wsComp.Range("C6").Value = 70.03886086508541
$wsComp.Range("C7").Value = 29.9611391349146
